$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.478.81'
$ws.Range("E2").Value = '  +4.96%  '
$ws.Range("D3").Value = '2.052.51'
$ws.Range("E3").Value = '  +3.37%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.79'
$ws.Range("E5").Value = '  +3.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.96'
$ws.Range("E7").Value = '  +10.84%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.400'
$ws.Range("E9").Value = '  +9.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '60.12'
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("E11").Value = '  +6.03%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.932'
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.77'
$ws.Range("E14").Value = '  +26.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.98'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").Value = '2.354.86'
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.77'
$ws.Range("E17").Value = '  +8.33%  '
$ws.Range("D18").Value = '2.062.07'
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").Value = '37.379.30'
$ws.Range("E19").Value = '  +4.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.92'
$ws.Range("E20").Value = '  +3.04%  '
$ws.Range("D21").Value = '0.0₃0893'
$ws.Range("E21").Value = '  +4.79%  '
$ws.Range("E22").Value = '  +5.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.24'
$ws.Range("E23").Value = '  +3.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.67'
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  +5.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.09'
$ws.Range("E27").Value = '  +9.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.26'
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.16'
$ws.Range("E29").Value = '  +4.46%  '
$ws.Range("E30").Value = '  +39.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.123'
$ws.Range("E31").Value = '  +3.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.25'
$ws.Range("E32").Value = '  +6.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("E33").Value = '  +6.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0636'
$ws.Range("E34").Value = '  +6.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.73'
$ws.Range("E35").Value = '  +8.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("E37").Value = '  +14.15%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.83'
$ws.Range("E39").Value = '  +2.60%  '
$ws.Range("E40").Value = '  +34.09%  '
$ws.Range("E41").Value = '  +11.26%  '
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("E43").Value = '  +7.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.84'
$ws.Range("E44").Value = '  +8.83%  '
$ws.Range("E45").Value = '  +7.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0221'
$ws.Range("E46").Value = '  +3.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.28'
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.00'
$ws.Range("D49").Value = '1.415.27'
$ws.Range("E49").Value = '  +3.72%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.95'
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '48.46'
$ws.Range("E51").Value = '  +2.62%  '
